# Modulo_1/M1_Objetivo_Ciencia_De_Dados.xlsx
#
# - "Registro Diário" (sheet 1): E3's method changes from
#   "Vídeo - Leitura" to the already-existing "Vídeo - Leitura - Pratica"
#   string, which makes "Vídeo - Leitura" unused and drops it from the
#   shared-strings table on save. The sheet also stops being the active tab
#   and its selection moves to E4.
# - "Progressso Geral" (sheet 3): gains a new row (A3 = 2), becomes the
#   active tab/sheet, and its selection moves to C4.

$wb = $excel.ActiveWorkbook

$registro = $wb.Worksheets.Item("Registro Diário")
$progresso = $wb.Worksheets.Item("Progressso Geral")

# Registro Diário: update the method used for the 2nd study entry.
$registro.Range("E3").Value = "Vídeo - Leitura - Pratica"
$registro.Range("E4").Select()

# Progressso Geral: log another completed module and make it the active tab.
$progresso.Range("A3").Value = 2
$progresso.Activate()
$progresso.Range("C4").Select()
